$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "coursesInterested" column (column D) entirely, shifting
# designation/collegeName left into D/E.
$ws.Range("D1").EntireColumn.Delete()
